# Handle mistakenly resolved checks
# Rows 85-91 were marked as voided/resolved in error; clear the check
# date (column A) and check number (column C) for those rows while
# leaving the row formatting untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A85:A91").ClearContents()
$ws.Range("C85:C91").ClearContents()

# Move the active selection up to where the last populated row now is.
$ws.Range("D91").Select()
